# Fruta / hortaliza, semanal
# Insert a new weekly observation row into the "Plátano" (Feria Lagunitas de
# Puerto Montt) consolidated dataset. The new record is inserted immediately
# above the existing row 963, pushing all subsequent rows down by one; the
# sheet's used range grows from A1:T1015 to A1:T1016.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 963 (shifts rows 963:1015 down to 964:1016).
$ws.Rows.Item(963).Insert()

# Populate the newly-inserted row with the new weekly record. Column layout
# mirrors every other row in the sheet:
# A Mercado ID | B Mercado | C Región | D Fecha | E Codreg | F Tipo
# G Producto ID | H Producto | I Categoría ID | J Categoría | K Variedad
# L Calidad | M Volumen | N Precio mínimo | O Precio máximo
# P Precio promedio ponderado | Q Unidad de comercialización | R Origen
# S Precio $/Kg | T Kg / unidad
$ws.Range("A963").Value2 = 4
$ws.Range("B963").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C963").Value2 = "Los Lagos"
$ws.Range("D963").Value2 = 45267
$ws.Range("E963").Value2 = 10
$ws.Range("F963").Value2 = "Fruta"
$ws.Range("G963").Value2 = 100108
$ws.Range("H963").Value2 = "Tropicales y subtropicales"
$ws.Range("I963").Value2 = 100108006
$ws.Range("J963").Value2 = "Plátano"
$ws.Range("K963").Value2 = "Sin especificar"
$ws.Range("L963").Value2 = "Primera Pintón"
$ws.Range("M963").Value2 = 800
$ws.Range("N963").Value2 = 29000
$ws.Range("O963").Value2 = 30000
$ws.Range("P963").Value2 = 29500
$ws.Range("Q963").Value2 = "$/caja 20 kilos"
$ws.Range("R963").Value2 = "Ecuador"
$ws.Range("S963").Value2 = 1475
$ws.Range("T963").Value2 = 20

# Date column keeps the sheet's date number format (style index 2).
$ws.Range("D963").NumberFormat = $ws.Range("D964").NumberFormat
